$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.353260278701782
$ws.Range("B1").Value = 2.559735298156738
$ws.Range("C1").Value = 2.013429403305054
$ws.Range("D1").Value = 1.902454018592834
$ws.Range("E1").Value = 1.699438452720642
